# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to match the refreshed cryptos list as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D updates that are plain numeric-looking strings: ---
# Force text format first so Excel keeps them as text cells (matches
# the original workbook where every Price cell is stored as a string),
# instead of silently re-interpreting them as numbers.
$numericPriceCells = @(
    @{ Cell = "D5"; Value = "213.03" }
    @{ Cell = "D6"; Value = "0.498" }
    @{ Cell = "D10"; Value = "18.99" }
    @{ Cell = "D11"; Value = "0.0850" }
    @{ Cell = "D16"; Value = "63.83" }
    @{ Cell = "D19"; Value = "215.28" }
    @{ Cell = "D20"; Value = "7.34" }
    @{ Cell = "D22"; Value = "4.29" }
    @{ Cell = "D23"; Value = "9.03" }
    @{ Cell = "D25"; Value = "144.86" }
    @{ Cell = "D27"; Value = "6.96" }
    @{ Cell = "D29"; Value = "15.10" }
    @{ Cell = "D30"; Value = "0.0490" }
    @{ Cell = "D39"; Value = "0.823" }
    @{ Cell = "D42"; Value = "0.942" }
    @{ Cell = "D44"; Value = "0.762" }
    @{ Cell = "D46"; Value = "60.88" }
    @{ Cell = "D47"; Value = "87.41" }
    @{ Cell = "D48"; Value = "1.48" }
    @{ Cell = "D50"; Value = "0.0953" }
)
foreach ($item in $numericPriceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# --- Column D updates that are already non-numeric text (safe as-is): ---
$ws.Range("D2").Value = "26.269.10"
$ws.Range("D3").Value = "1.593.47"
$ws.Range("D12").Value = "1.817.25"
$ws.Range("D13").Value = "1.581.19"
$ws.Range("D17").Value = "26.257.70"
$ws.Range("D18").Value = "0.0₃0721"
$ws.Range("D33").Value = "1.414.96"
$ws.Range("D45").Value = "1.729.26"

# --- Column E (Volume 1h %) updates: ---
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  +5.97%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("E37").Value = "  -5.26%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -11.24%  "
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("E51").Value = "  +0.07%  "

